# Apply edit: fill in the "Definition" column (D) on the "Concepts" sheet
# with the same text as the "Display" column (C) for each concept row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Copy the Display (C) value into the empty Definition (D) cell for each
# data row (rows 2-5).
for ($row = 2; $row -le 5; $row++) {
    $display = $ws.Cells.Item($row, 3).Text
    $ws.Cells.Item($row, 4).Value = $display
}
